# Update policy indicators to CDC's 11/10/2022 suggestions.
# Removes the old continuous/quartile-based dummy columns and replaces them
# with the new CDC-suggestion dummy indicator columns (C/D/E) plus updated
# Construct labels (column B) for a few rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the CDC-suggestion columns
$ws.Range("B1").Value = 'Construct'
$ws.Range("C1").Value = 'cdc'
$ws.Range("D1").Value = 'operation'
$ws.Range("E1").Value = 'score_needed'

    $ws.Range("B1").Value = 'Construct'
    $ws.Range("C1").Value = 'cdc'
    $ws.Range("D1").Value = 'operation'
    $ws.Range("E1").Value = 'score_needed'
    $ws.Range("C2").Value = 0
    $ws.Range("D2").Value = 'none'
    $ws.Range("C3").Value = 0
    $ws.Range("D3").Value = 'none'
    $ws.Range("C4").Value = 0
    $ws.Range("D4").Value = 'none'
    $ws.Range("C5").Value = 1
    $ws.Range("D5").Value = 'either'
    $ws.Range("E5").Value = 2
    $ws.Range("C6").Value = 1
    $ws.Range("D6").Value = 'either'
    $ws.Range("E6").Value = 2
    $ws.Range("C7").Value = 0
    $ws.Range("D7").Value = 'none'
    $ws.Range("C8").Value = 0
    $ws.Range("D8").Value = 'none'
    $ws.Range("C9").Value = 0
    $ws.Range("D9").Value = 'none'
    $ws.Range("C10").Value = 0
    $ws.Range("D10").Value = 'none'
    $ws.Range("C11").Value = 0
    $ws.Range("D11").Value = 'none'
    $ws.Range("C12").Value = 0
    $ws.Range("D12").Value = 'none'
    $ws.Range("C13").Value = 0
    $ws.Range("D13").Value = 'none'
    $ws.Range("C14").Value = 0
    $ws.Range("D14").Value = 'none'
    $ws.Range("C15").Value = 1
    $ws.Range("D15").Value = 'all'
    $ws.Range("E15").Value = 3
    $ws.Range("C16").Value = 1
    $ws.Range("D16").Value = 'all'
    $ws.Range("E16").Value = 3
    $ws.Range("C17").Value = 0
    $ws.Range("D17").Value = 'none'
    $ws.Range("C18").Value = 0
    $ws.Range("D18").Value = 'none'
    $ws.Range("C19").Value = 0
    $ws.Range("D19").Value = 'none'
    $ws.Range("C20").Value = 1
    $ws.Range("D20").Value = 'all'
    $ws.Range("E20").Value = 3
    $ws.Range("C21").Value = 0
    $ws.Range("D21").Value = 'none'
    $ws.Range("C22").Value = 0
    $ws.Range("D22").Value = 'none'
    $ws.Range("C23").Value = 0
    $ws.Range("D23").Value = 'none'
    $ws.Range("B24").Value = 'Screening testing for students'
    $ws.Range("C24").Value = 0
    $ws.Range("D24").Value = 'none'
    $ws.Range("C25").Value = 0
    $ws.Range("D25").Value = 'none'
    $ws.Range("B26").Value = 'Screening testing for students'
    $ws.Range("C26").Value = 1
    $ws.Range("D26").Value = 'all'
    $ws.Range("E26").Value = 2
    $ws.Range("C27").Value = 0
    $ws.Range("D27").Value = 'none'
    $ws.Range("C28").Value = 0
    $ws.Range("D28").Value = 'none'
    $ws.Range("C29").Value = 0
    $ws.Range("D29").Value = 'none'
    $ws.Range("C30").Value = 1
    $ws.Range("D30").Value = 'all'
    $ws.Range("E30").Value = 3
    $ws.Range("B31").Value = 'Contact tracing'
    $ws.Range("C31").Value = 1
    $ws.Range("D31").Value = 'all'
    $ws.Range("E31").Value = 2
    $ws.Range("B32").Value = 'Contact tracing'
    $ws.Range("C32").Value = 0
    $ws.Range("D32").Value = 'none'
    $ws.Range("B33").Value = 'Quarantine'
    $ws.Range("C33").Value = 0
    $ws.Range("D33").Value = 'none'
    $ws.Range("B34").Value = 'Quarantine'
    $ws.Range("C34").Value = 1
    $ws.Range("D34").Value = 'all'
    $ws.Range("E34").Value = 3
    $ws.Range("C35").Value = 0
    $ws.Range("D35").Value = 'none'
    $ws.Range("C36").Value = 1
    $ws.Range("D36").Value = 'all'
    $ws.Range("E36").Value = 3
    $ws.Range("B37").Value = 'Ventilation'
    $ws.Range("C37").Value = 0
    $ws.Range("D37").Value = 'none'
    $ws.Range("B38").Value = 'Ventilation'
    $ws.Range("C38").Value = 0
    $ws.Range("D38").Value = 'none'
    $ws.Range("B39").Value = 'Ventilation'
    $ws.Range("C39").Value = 0
    $ws.Range("D39").Value = 'none'
    $ws.Range("B40").Value = 'Ventilation'
    $ws.Range("C40").Value = 0
    $ws.Range("D40").Value = 'none'
    $ws.Range("B41").Value = 'HEPA filters'
    $ws.Range("C41").Value = 1
    $ws.Range("D41").Value = 'all'
    $ws.Range("E41").Value = 2
    $ws.Range("B42").Value = 'HVAC systems'
    $ws.Range("C42").Value = 1
    $ws.Range("D42").Value = 'all'
    $ws.Range("E42").Value = 2

# Column B ("Construct") needs to be a bit wider now that it holds the new
# CDC construct labels.
$ws.Columns("B").ColumnWidth = 20.6640625

# Restore the active cell/selection used when the workbook was last saved.
$ws.Range("B25").Select()

